## Persist null members in array (#28)
## Adds a new "BadSheet" worksheet (after Sheet2) that demonstrates
## range evaluation fault tolerance around holes (blank cells) inside
## a lookup range, plus a bump of the volatile NOW() cache + a window
## view tweak on Sheet1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1: workbook/window view tweaks + refreshed NOW() value
# ---------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$excel.ActiveWindow.ScrollRow = 1
$sheet1.Range("B25").Select()
$sheet1.Range("B2").Calculate()

# ---------------------------------------------------------------
# New worksheet "BadSheet", inserted right after Sheet2
# ---------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$ws = $wb.Worksheets.Add($null, $sheet2)
$ws.Name = "BadSheet"

$ws.Columns.Item(1).ColumnWidth = 15.21875

# Title banner, merged across A1:H1, bold red left-aligned text.
$title = $ws.Range("A1:H1")
$title.Value = "README: THIS IS A BAD SHEET TO DEMONSTRATE FAULT TOLERANCE AND ERROR RESISTANCE"
$title.Font.Bold = $true
$title.Font.Color = 255
$title.HorizontalAlignment = -4131
$title.Merge()

# "Range with holes" section header (A3:D3) + input grid (B4:D6) with
# a couple of intentionally-empty "holes" at B4 and C6.
$ws.Range("A3:D3").Value = "Range with holes"
$ws.Range("A3").Style = "Normal"

$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 19
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 34
$ws.Range("B6").Value = 3
$ws.Range("D6").Value = 39

# "RANGE EVALUATION" section header (A8:D8)
$ws.Range("A8:D8").Value = "RANGE EVALUATION"

# FORMULATEXT() echoes alongside SUM / VLOOKUP probes against the
# holey range above - rows 9-12.
$ws.Range("A9").Formula = '=FORMULATEXT(B9)'
$ws.Range("B9").Formula = '=SUM(B4:C6)'

$ws.Range("A10").Formula = '=FORMULATEXT(B10)'
$ws.Range("B10").Formula = '=VLOOKUP(0,$B$4:$C$6,2,0)'

$ws.Range("A11").Formula = '=FORMULATEXT(B11)'
$ws.Range("B11").Formula = '=VLOOKUP(2,$B$4:$C$6,2,0)'

$ws.Range("A12").Formula = '=FORMULATEXT(B12)'
$ws.Range("B12").Formula = '=VLOOKUP(3,$B$4:$C$6,2,0)'

$ws.Range("C17").Select()
